# Add team record (Wins/Losses/Ties) columns to the CLE_2019 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD="Wins", AE="Losses", AF="Ties" ---
# Copy the formatting from the existing header cell (AC1) so the new
# headers share the same bold/centered/bordered style used by the rest
# of row 1, then set the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-56: team record for every player (93-69-0) ---
$lastRow = 56
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 93  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 69  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
